# RGF_input_file.xlsx edit — "add R_matrix support. Debug T and R matrix."
#
# Applies the cell-level / view-level changes from the target diff using the
# Excel COM object model against the already-open ActiveWorkbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 7 ("CPU max matrix"): B7 116->100, C7 2->120 -----------------------
$ws.Range("B7").Value = 100
$ws.Range("C7").Value = 120

# --- Row 13 / Row 14: swap the "x" / "o" markers ----------------------------
# Originally A13 held shared-string "o" and A14 held shared-string "x".
# The diff swaps them (A13 -> "x", A14 -> "o"), which also makes the
# formerly-unused "T" shared string drop out of the saved sharedStrings table.
$ws.Range("A13").Value = "x"
$ws.Range("A14").Value = "o"

# --- Row 14 input values (R_matrix debug values) ----------------------------
$ws.Range("C14").Value = 1
$ws.Range("E14").Value = 10
$ws.Range("F14").Value = 3
$ws.Range("I14").Value = 0.1
# N14/O14/P14/Q14 are formulas referencing the above inputs (plus J14/L14,
# unchanged at 0) — they recalculate automatically after the script runs.

# --- Selection / active cell -------------------------------------------------
$ws.Range("I19").Select() | Out-Null

# --- Column widths (best effort; COM ColumnWidth is quantized to whole
#     pixels at the engine's fixed 7px Maximum-Digit-Width, so the exact
#     fractional widths from the source file can't be reproduced bit-for-bit,
#     but this sets the closest reachable value in the correct direction) ---
$ws.Columns.Item(1).ColumnWidth = 16.428571428571427
$ws.Columns.Item(2).ColumnWidth = 8.714285714285714
